$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados" timestamp bumped from 22:35 to 23:05 ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 23:05"

# --- Countries list reorder: "Somalia" moves up to right after "Nueva Zelanda",
#     pushing "El Salvador", "Eslovaquia" and "Eslovenia" down one row each.
#     ("Gabon" on row 98 is unaffected.)
$ws.Range("A94").Value = "Somalia"
$ws.Range("A95").Value = "El Salvador"
$ws.Range("A96").Value = "Eslovaquia"
$ws.Range("A97").Value = "Eslovenia"

# --- Updated case counters ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1565229
$ws.Range("C4").Value = 14935
$ws.Range("D4").Value = 361854
$ws.Range("E4").Value = 1110218
$ws.Range("G4").Value = 1176
$ws.Range("H4").Value = 93157

# Row 11
$ws.Range("B11").Value = 177824
$ws.Range("C11").Value = 535
$ws.Range("E11").Value = 13931
$ws.Range("G11").Value = 70
$ws.Range("H11").Value = 8193

# Row 14
$ws.Range("B14").Value = 106475
$ws.Range("C14").Value = 6147
$ws.Range("D14").Value = 42309
$ws.Range("E14").Value = 60864
$ws.Range("G14").Value = 146
$ws.Range("H14").Value = 3302

# Row 83
$ws.Range("B83").Value = 2153
$ws.Range("C83").Value = 34
$ws.Range("D83").Value = 1050
$ws.Range("E83").Value = 1075

# Row 94: now "Somalia" - new/updated statistics
$ws.Range("B94").Value = 1502
$ws.Range("C94").Value = 47
$ws.Range("D94").Value = 178
$ws.Range("E94").Value = 1265
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 59

# Row 95: now "El Salvador" (former row-94 statistics, unchanged values, shifted down)
$ws.Range("B95").Value = 1498
$ws.Range("C95").Value = 85
$ws.Range("D95").Value = 502
$ws.Range("E95").Value = 966
$ws.Range("H95").Value = 30

# Row 96: now "Eslovaquia" (former row-95 statistics, unchanged values, shifted down)
$ws.Range("B96").Value = 1495
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 1192
$ws.Range("E96").Value = 275
$ws.Range("H96").Value = 28

# Row 97: now "Eslovenia" (former row-96 statistics, unchanged values, shifted down)
$ws.Range("B97").Value = 1467
$ws.Range("C97").Value = 1
$ws.Range("D97").Value = 1335
$ws.Range("E97").Value = 28
$ws.Range("H97").Value = 104

# Row 123
$ws.Range("D123").Value = 417
$ws.Range("E123").Value = 223
